$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 39
$ws.Range("B2").Value = "5:45 AM"
$ws.Range("C2").Value = 1169.48
$ws.Range("D2").Value = 1199.48
$ws.Range("E2").Value = 21.87

$ws.Range("A3").Value = 49
$ws.Range("B3").Value = "5:50 AM"
$ws.Range("C3").Value = 1467.18
$ws.Range("D3").Value = 1497.18
$ws.Range("E3").Value = 24.93

$ws.Range("A4").Value = 91
$ws.Range("B4").Value = "6:11 AM"
$ws.Range("C4").Value = 2725.88
$ws.Range("D4").Value = 2755.88
$ws.Range("E4").Value = 22.16

$ws.Range("A5").Value = 100
$ws.Range("B5").Value = "6:15 AM"
$ws.Range("C5").Value = 2986.3175
$ws.Range("D5").Value = 3016.3175
$ws.Range("E5").Value = 30.9975

$ws.Range("A6").Value = 120
$ws.Range("B6").Value = "6:25 AM"
$ws.Range("C6").Value = 3577.16
$ws.Range("D6").Value = 3607.16
$ws.Range("E6").Value = 49.37

$ws.Range("A8").Value = 239
$ws.Range("B8").Value = "7:24 AM"
$ws.Range("C8").Value = 7146.62
$ws.Range("D8").Value = 7176.62
$ws.Range("E8").Value = 32.95

$ws.Range("A9").Value = 257
$ws.Range("B9").Value = "7:33 AM"
$ws.Range("C9").Value = 7691.65
$ws.Range("D9").Value = 7721.65
$ws.Range("E9").Value = 35.775

$ws.Range("A10").Value = 275
$ws.Range("B10").Value = "7:43 AM"
$ws.Range("C10").Value = 8248.530000000001
$ws.Range("D10").Value = 8278.530000000001
$ws.Range("E10").Value = 21.35

$ws.Range("A13").Value = 570
$ws.Range("B13").Value = "10:10 A"
$ws.Range("C13").Value = 17087.03
$ws.Range("D13").Value = 17117.03
$ws.Range("E13").Value = 37.05

$ws.Range("A14").Value = 993
$ws.Range("B14").Value = "1:41 PM"
$ws.Range("C14").Value = 29775.52
$ws.Range("D14").Value = 29805.52
$ws.Range("E14").Value = 41.35

$ws.Range("A15").Value = 1016
$ws.Range("B15").Value = "1:53 PM"
$ws.Range("C15").Value = 30454.73
$ws.Range("D15").Value = 30484.73
$ws.Range("E15").Value = 24.88

$ws.Range("A16").Value = 1022
$ws.Range("B16").Value = "1:56 PM"
$ws.Range("C16").Value = 30644.296667
$ws.Range("D16").Value = 30674.296667
$ws.Range("E16").Value = 24.493333

$ws.Range("A17").Value = 1036
$ws.Range("B17").Value = "2:03 PM"
$ws.Range("C17").Value = 31058.49
$ws.Range("D17").Value = 31088.49
$ws.Range("E17").Value = 60.24

$ws.Range("A18").Value = 1086
$ws.Range("B18").Value = "2:28 PM"
$ws.Range("C18").Value = 32575.77
$ws.Range("D18").Value = 32605.77
$ws.Range("E18").Value = 21.33

$ws.Range("A19").Value = 1103
$ws.Range("B19").Value = "2:36 PM"
$ws.Range("C19").Value = 33070.88
$ws.Range("D19").Value = 33100.88
$ws.Range("E19").Value = 25.1

$ws.Range("A20").Value = 1141
$ws.Range("B20").Value = "2:55 PM"
$ws.Range("C20").Value = 34219.87
$ws.Range("D20").Value = 34249.87
$ws.Range("E20").Value = 21.48

$ws.Range("A21").Value = 1185
$ws.Range("B21").Value = "3:17 PM"
$ws.Range("C21").Value = 35530.52
$ws.Range("D21").Value = 35560.52
$ws.Range("E21").Value = 115.925

$ws.Range("A22").Value = 1194
$ws.Range("B22").Value = "3:22 PM"
$ws.Range("C22").Value = 35803.485
$ws.Range("D22").Value = 35833.485
$ws.Range("E22").Value = 24.31

$ws.Range("A23").Value = 1202
$ws.Range("B23").Value = "3:26 PM"
$ws.Range("C23").Value = 36041.01
$ws.Range("D23").Value = 36071.01
$ws.Range("E23").Value = 24.126667

$ws.Range("A24").Value = 1209
$ws.Range("B24").Value = "3:29 PM"
$ws.Range("C24").Value = 36244.02
$ws.Range("D24").Value = 36274.02
$ws.Range("E24").Value = 75.29000000000001

$ws.Range("A25").Value = 1221
$ws.Range("B25").Value = "3:35 PM"
$ws.Range("C25").Value = 36616.4
$ws.Range("D25").Value = 36646.4
$ws.Range("E25").Value = 23.47

$ws.Range("A26").Value = 1227
$ws.Range("B26").Value = "3:39 PM"
$ws.Range("C26").Value = 36808.87
$ws.Range("D26").Value = 36838.87
$ws.Range("E26").Value = 22.12

$ws.Range("A27").Value = 1245
$ws.Range("B27").Value = "3:47 PM"
$ws.Range("C27").Value = 37332.09
$ws.Range("D27").Value = 37362.09
$ws.Range("E27").Value = 193.12

$ws.Range("A28").Value = 1257
$ws.Range("B28").Value = "3:53 PM"
$ws.Range("C28").Value = 37684.8
$ws.Range("D28").Value = 37714.8
$ws.Range("E28").Value = 36.82

$ws.Range("A29").Value = 1269
$ws.Range("B29").Value = "3:59 PM"
$ws.Range("C29").Value = 38052.733333
$ws.Range("D29").Value = 38082.733333
$ws.Range("E29").Value = 20.703333

$ws.Range("A33").Value = 1308
$ws.Range("B33").Value = "4:19 PM"
$ws.Range("C33").Value = 39217.913333
$ws.Range("D33").Value = 39247.913333
$ws.Range("E33").Value = 29.03

$ws.Range("A34").Value = 1318
$ws.Range("B34").Value = "4:24 PM"
$ws.Range("C34").Value = 39526.18
$ws.Range("D34").Value = 39556.18
$ws.Range("E34").Value = 111.31

$ws.Range("A35").Value = 1335
$ws.Range("B35").Value = "4:32 PM"
$ws.Range("C35").Value = 40027.72
$ws.Range("D35").Value = 40057.72
$ws.Range("E35").Value = 47.705

$ws.Range("A36").Value = 1341
$ws.Range("B36").Value = "4:35 PM"
$ws.Range("C36").Value = 40209.895
$ws.Range("D36").Value = 40239.895
$ws.Range("E36").Value = 54.99

$ws.Range("A37").Value = 1354
$ws.Range("B37").Value = "4:42 PM"
$ws.Range("C37").Value = 40599.6
$ws.Range("D37").Value = 40629.6
$ws.Range("E37").Value = 20.99

$ws.Range("A38").Value = 1362
$ws.Range("B38").Value = "4:46 PM"
$ws.Range("C38").Value = 40838.935
$ws.Range("D38").Value = 40868.935
$ws.Range("E38").Value = 104.33

$ws.Range("A39").Value = 1372
$ws.Range("B39").Value = "4:51 PM"
$ws.Range("C39").Value = 41141.465
$ws.Range("D39").Value = 41171.465
$ws.Range("E39").Value = 98.08

$ws.Range("A41").Value = 1410
$ws.Range("B41").Value = "5:10 PM"
$ws.Range("C41").Value = 42278.335
$ws.Range("D41").Value = 42308.335
$ws.Range("E41").Value = 44.55

$ws.Range("A42").Value = 1417
$ws.Range("B42").Value = "5:13 PM"
$ws.Range("C42").Value = 42486.06
$ws.Range("D42").Value = 42516.06
$ws.Range("E42").Value = 61.41

$ws.Range("A43").Value = 1432
$ws.Range("B43").Value = "5:21 PM"
$ws.Range("C43").Value = 42943.82
$ws.Range("D43").Value = 42973.82
$ws.Range("E43").Value = 31.495

$ws.Range("A51").Value = 1546
$ws.Range("B51").Value = "6:18 PM"
$ws.Range("C51").Value = 46356.505
$ws.Range("D51").Value = 46386.505
$ws.Range("E51").Value = 25.355
